# "materi 2 vlookup excel midle selesai"
# Add a VLOOKUP "Cek Stok" lookup column (F) on sheet4, driven by the
# fruit names already typed into column E, and extend the E/F sample data
# down to row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet4")

# F2:F4 already have their lookup value in column E - just add the formula.
$ws.Range("F2").Formula = '=IFERROR(VLOOKUP(E2,A1:$C$6,3,0),"No Data")'
$ws.Range("F3").Formula = '=IFERROR(VLOOKUP(E3,A2:$C$6,3,0),"No Data")'
$ws.Range("F4").Formula = '=IFERROR(VLOOKUP(E4,A3:$C$6,3,0),"No Data")'

# Row 5 gets a new lookup value ("Melon") plus its VLOOKUP formula.
$ws.Range("E5").Value = "Melon"
$ws.Range("F5").Formula = '=IFERROR(VLOOKUP(E5,A4:$C$6,3,0),"No Data")'

# Leave the new formula block selected, like the author did after typing it.
$ws.Range("F2:F5").Select() | Out-Null
